# GPLIM-5135 fixes and improvements for Pooled Tube Upload
# Populate the "Read Length" value for the sample row and leave the
# selection on the cell that was just filled in, matching the
# round-trip produced by Excel after the data entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 61
$ws.Range("P2").Select() | Out-Null
